$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade the "Practical" assignment (new comment/grade pair in columns D/E) for the
# students who submitted something, row by row first...

$ws.Range("D2").Value = "Good job! Next time please comment answers"
$ws.Range("E2").Value = 6

$ws.Range("D3").Value = "Great job!"
$ws.Range("E3").Value = 6

$ws.Range("D4").Value = "Great job!"
$ws.Range("E4").Value = 6

$ws.Range("D5").Value = "Good job! Unable to do a few using R but left everything commented"
$ws.Range("E5").Value = 5.75

$ws.Range("D6").Value = "Good job! Updated 2 files but all seems ok"
$ws.Range("E6").Value = 6

$ws.Range("D7").Value = "Great job!"
$ws.Range("E7").Value = 6

$ws.Range("D9").Value = "Great job!"
$ws.Range("E9").Value = 6

$ws.Range("D10").Value = "Sent us an assignment with 3 questions only?"
$ws.Range("E10").Value = 3

$ws.Range("D11").Value = "Great job!"
$ws.Range("E11").Value = 6

$ws.Range("D12").Value = "Great job!"
$ws.Range("E12").Value = 6

$ws.Range("D13").Value = "Good job! Not the most efficient code but got things done"
$ws.Range("E13").Value = 6

$ws.Range("D14").Value = "Great job!"
$ws.Range("E14").Value = 6

$ws.Range("D15").Value = "Great job!"
$ws.Range("E15").Value = 6

# ...then go back and mark the students who never submitted anything.
$ws.Range("D8").Value = "No submission - No justification"

$ws.Range("D16").Value = "Asked for an extension beforehand…"
$ws.Range("E16").Value = 0

$ws.Range("B17").Value = "No submission - No justification"
$ws.Range("D17").Value = "No submission - No justification"
$ws.Range("E17").Value = 0

$ws.Range("B18").Value = "No submission - No justification"
$ws.Range("D18").Value = "No submission - No justification"
$ws.Range("E18").Value = 0

$ws.Range("B19").Value = "No submission - No justification"
$ws.Range("D19").Value = "No submission - No justification"
$ws.Range("E19").Value = 0

$ws.Range("B20").Value = "No submission - No justification"
$ws.Range("D20").Value = "No submission - No justification"
$ws.Range("E20").Value = 0

$ws.Range("B21").Value = "No submission - No justification"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = "No submission - No justification"
$ws.Range("E21").Value = 0

# Column D now holds the grading comments - widen it to match the other comment columns
$ws.Columns.Item(4).ColumnWidth = 61.83203125

# Leave the selection on the last cell that was graded
[void]$ws.Range("A21").Select()
